# Update cryptocurrency price/volume data per latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.325.13"
$ws.Range("E2").Value = '  +3.06%  '

$ws.Range("D3").Value = "'2.679.27"
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").Value = "'523.02"
$ws.Range("E5").Value = '  +2.17%  '

$ws.Range("D6").Value = "'145.95"
$ws.Range("E6").Value = '  +2.28%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("E8").Value = '  +2.15%  '

$ws.Range("D9").Value = "'2.699.31"
$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("E10").Value = '  +2.89%  '

$ws.Range("E11").Value = '  +0.60%  '

$ws.Range("E12").Value = '  +2.42%  '

$ws.Range("E13").Value = '  +1.82%  '

$ws.Range("D14").Value = "'3.152.29"
$ws.Range("E14").Value = '  +1.48%  '

$ws.Range("D15").Value = "'60.374.83"
$ws.Range("E15").Value = '  +3.12%  '

$ws.Range("D16").Value = "'21.26"
$ws.Range("E16").Value = '  +1.97%  '

$ws.Range("D17").Value = "'2.763.37"
$ws.Range("E17").Value = '  +3.55%  '

$ws.Range("E18").Value = '  +1.87%  '

$ws.Range("D19").Value = "'350.44"
$ws.Range("E19").Value = '  +2.81%  '

$ws.Range("D20").Value = "'4.54"
$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").Value = "'10.61"
$ws.Range("E21").Value = '  +2.33%  '

$ws.Range("E22").Value = '  +3.80%  '

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").Value = "'62.76"
$ws.Range("E24").Value = '  +3.22%  '

$ws.Range("D25").Value = "'0.423"
$ws.Range("E25").Value = '  +1.21%  '

$ws.Range("E26").Value = '  +5.85%  '

$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("D28").Value = "'0.0₃0815"
$ws.Range("E28").Value = '  +1.70%  '

$ws.Range("D29").Value = "'7.25"
$ws.Range("E29").Value = '  +1.50%  '

$ws.Range("D30").Value = "'6.84"
$ws.Range("E30").Value = '  +7.98%  '

$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("E32").Value = '  +1.91%  '

$ws.Range("D33").Value = "'19.09"
$ws.Range("E33").Value = '  +1.36%  '

$ws.Range("D34").Value = "'148.01"
$ws.Range("E34").Value = '  -0.94%  '

$ws.Range("D35").Value = "'4.31"
$ws.Range("E35").Value = '  +8.25%  '

$ws.Range("E36").Value = '  +8.99%  '

$ws.Range("D37").Value = "'0.952"
$ws.Range("E37").Value = '  -5.06%  '

$ws.Range("D38").Value = "'1.54"
$ws.Range("E38").Value = '  +11.09%  '

$ws.Range("D39").Value = "'0.877"
$ws.Range("E39").Value = '  +3.42%  '

$ws.Range("D40").Value = "'36.85"
$ws.Range("E40").Value = '  +0.97%  '

$ws.Range("E41").Value = '  +0.96%  '

$ws.Range("D42").Value = "'281.41"
$ws.Range("E42").Value = '  +1.13%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = "'0.997"
$ws.Range("E43").Value = '  +0.29%  '

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = "'0.0988"
$ws.Range("E44").Value = '  +1.36%  '

$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = "'0.611"
$ws.Range("E45").Value = '  -0.80%  '

$ws.Range("E46").Value = '  +2.25%  '

$ws.Range("D47").Value = "'2.137.15"
$ws.Range("E47").Value = '  +7.68%  '

$ws.Range("D48").Value = "'0.0538"
$ws.Range("E48").Value = '  +1.53%  '

$ws.Range("D49").Value = "'4.87"
$ws.Range("E49").Value = '  +3.95%  '

$ws.Range("D50").Value = "'0.0234"
$ws.Range("E50").Value = '  +2.47%  '

$ws.Range("D51").Value = "'10.45"
$ws.Range("E51").Value = '  +1.85%  '
